$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.0
$ws.Range("C2").Value = "dsfj"
$ws.Range("D2").Value = "'234"
$ws.Range("E2").Value = "kf"
